# Refresh the FlashScore weekly fixtures sheet:
#  - insert a new "Odd_CS_4-4" header column (shifting the CS_0-1.. CS_2-3 block right by one)
#  - insert a new match row (GV San Jose vs Aurora) as row 2
#  - keep the two existing matches, re-aligned to the new column order, as rows 3-4
#  - append a new match row (Eldense vs Huesca) as row 5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: a new "Odd_CS_4-4" column is inserted before "Odd_CS_0-1" (was at the end of the CS block, now at AG); the rest of the CS_x-y block shifts right by one column.
$row1 = New-Object 'object[,]' 1,56
$row1[0,0] = "Id"
$row1[0,1] = "Date"
$row1[0,2] = "Time"
$row1[0,3] = "League"
$row1[0,4] = "Home"
$row1[0,5] = "Away"
$row1[0,6] = "Odd_H_FT"
$row1[0,7] = "Odd_D_FT"
$row1[0,8] = "Odd_A_FT"
$row1[0,9] = "Odd_H_HT"
$row1[0,10] = "Odd_D_HT"
$row1[0,11] = "Odd_A_HT"
$row1[0,12] = "Odd_Over05_FT"
$row1[0,13] = "Odd_Under05_FT"
$row1[0,14] = "Odd_Over15_FT"
$row1[0,15] = "Odd_Under15_FT"
$row1[0,16] = "Odd_Over25_FT"
$row1[0,17] = "Odd_Under25_FT"
$row1[0,18] = "Odd_Over05_HT"
$row1[0,19] = "Odd_Under05_HT"
$row1[0,20] = "Odd_BTTS_Yes"
$row1[0,21] = "Odd_BTTS_No"
$row1[0,22] = "Odd_CS_1-0"
$row1[0,23] = "Odd_CS_2-0"
$row1[0,24] = "Odd_CS_2-1"
$row1[0,25] = "Odd_CS_3-0"
$row1[0,26] = "Odd_CS_3-1"
$row1[0,27] = "Odd_CS_3-2"
$row1[0,28] = "Odd_CS_0-0"
$row1[0,29] = "Odd_CS_1-1"
$row1[0,30] = "Odd_CS_2-2"
$row1[0,31] = "Odd_CS_3-3"
$row1[0,32] = "Odd_CS_4-4"
$row1[0,33] = "Odd_CS_0-1"
$row1[0,34] = "Odd_CS_0-2"
$row1[0,35] = "Odd_CS_1-2"
$row1[0,36] = "Odd_CS_0-3"
$row1[0,37] = "Odd_CS_1-3"
$row1[0,38] = "Odd_CS_2-3"
$row1[0,39] = "Odd_CS_1-0_HT"
$row1[0,40] = "Odd_CS_2-0_HT"
$row1[0,41] = "Odd_CS_2-1_HT"
$row1[0,42] = "Odd_CS_3-0_HT"
$row1[0,43] = "Odd_CS_3-1_HT"
$row1[0,44] = "Odd_CS_3-2_HT"
$row1[0,45] = "Odd_CS_0-0_HT"
$row1[0,46] = "Odd_CS_1-1_HT"
$row1[0,47] = "Odd_CS_2-2_HT"
$row1[0,48] = "Odd_CS_0-1_HT"
$row1[0,49] = "Odd_CS_0-2_HT"
$row1[0,50] = "Odd_CS_1-2_HT"
$row1[0,51] = "Odd_CS_0-3_HT"
$row1[0,52] = "Odd_CS_1-3_HT"
$row1[0,53] = "Odd_CS_2-3_HT"
$row1[0,54] = "Odd_CS_3-3_HT"
$row1[0,55] = "Odd_CS_4-4_HT"
$ws.Range("A1:BD1").Value = $row1

# New row inserted for GV San Jose vs Aurora (Bolivia - Division Profesional, 16:00).
$row2 = New-Object 'object[,]' 1,56
$row2[0,0] = "GbLZ1izi"
$row2[0,1] = "20/11/2024"
$row2[0,2] = "16:00"
$row2[0,3] = "BOLIVIA - DIVISION PROFESIONAL"
$row2[0,4] = "GV San Jose"
$row2[0,5] = "Aurora"
$row2[0,6] = 1.9
$row2[0,7] = 3.7
$row2[0,8] = 3.7
$row2[0,9] = 2.5
$row2[0,10] = 2.3
$row2[0,11] = 4
$row2[0,12] = 1.03
$row2[0,13] = 15
$row2[0,14] = 1.2
$row2[0,15] = 4.33
$row2[0,16] = 1.67
$row2[0,17] = 2.15
$row2[0,18] = 1.33
$row2[0,19] = 3.25
$row2[0,20] = 1.62
$row2[0,21] = 2.2
$row2[0,22] = 9
$row2[0,23] = 10
$row2[0,24] = 8.5
$row2[0,25] = 17
$row2[0,26] = 15
$row2[0,27] = 21
$row2[0,28] = 13
$row2[0,29] = 7.5
$row2[0,30] = 13
$row2[0,31] = 41
$row2[0,32] = 151
$row2[0,33] = 13
$row2[0,34] = 21
$row2[0,35] = 13
$row2[0,36] = 41
$row2[0,37] = 26
$row2[0,38] = 29
$row2[0,39] = 4
$row2[0,40] = 10
$row2[0,41] = 19
$row2[0,42] = 34
$row2[0,43] = 51
$row2[0,44] = 101
$row2[0,45] = 3.25
$row2[0,46] = 7.5
$row2[0,47] = 41
$row2[0,48] = 6
$row2[0,49] = 19
$row2[0,50] = 23
$row2[0,51] = 51
$row2[0,52] = 67
$row2[0,53] = 151
$row2[0,54] = 51
$row2[0,55] = 51
$ws.Range("A2:BD2").Value = $row2

# Existing Corinthians vs Cruzeiro row (Brazil), now shifted down to row 3; CS odds columns AG:AM re-aligned to match the new header order.
$row3 = New-Object 'object[,]' 1,56
$row3[0,0] = "vgnPzklr"
$row3[0,1] = "20/11/2024"
$row3[0,2] = "11:00"
$row3[0,3] = "BRAZIL - SERIE A BETANO"
$row3[0,4] = "Corinthians"
$row3[0,5] = "Cruzeiro"
$row3[0,6] = 1.4
$row3[0,7] = 4.5
$row3[0,8] = 8
$row3[0,9] = 1.95
$row3[0,10] = 2.3
$row3[0,11] = 7.5
$row3[0,12] = 1.05
$row3[0,13] = 11
$row3[0,14] = 1.3
$row3[0,15] = 3.5
$row3[0,16] = 1.98
$row3[0,17] = 1.92
$row3[0,18] = 1.4
$row3[0,19] = 2.75
$row3[0,20] = 2.2
$row3[0,21] = 1.62
$row3[0,22] = 6
$row3[0,23] = 6
$row3[0,24] = 9
$row3[0,25] = 8.5
$row3[0,26] = 13
$row3[0,27] = 34
$row3[0,28] = 10
$row3[0,29] = 9
$row3[0,30] = 23
$row3[0,31] = 81
$row3[0,32] = 201
$row3[0,33] = 15
$row3[0,34] = 41
$row3[0,35] = 21
$row3[0,36] = 101
$row3[0,37] = 51
$row3[0,38] = 51
$row3[0,39] = 3.2
$row3[0,40] = 7
$row3[0,41] = 21
$row3[0,42] = 21
$row3[0,43] = 51
$row3[0,44] = 151
$row3[0,45] = 2.75
$row3[0,46] = 10
$row3[0,47] = 67
$row3[0,48] = 9
$row3[0,49] = 41
$row3[0,50] = 41
$row3[0,51] = 201
$row3[0,52] = 201
$row3[0,53] = 451
$row3[0,54] = 126
$row3[0,55] = 126
$ws.Range("A3:BD3").Value = $row3

# Existing Grobina vs Alberts JDFS row (Latvia), now shifted down to row 4; CS odds columns AG:AM re-aligned to match the new header order.
$row4 = New-Object 'object[,]' 1,56
$row4[0,0] = "468pA9I6"
$row4[0,1] = "20/11/2024"
$row4[0,2] = "08:00"
$row4[0,3] = "LATVIA - VIRSLIGA"
$row4[0,4] = "Grobina"
$row4[0,5] = "Alberts JDFS"
$row4[0,6] = 1.34
$row4[0,7] = 4.55
$row4[0,8] = 6.9
$row4[0,9] = 1.8
$row4[0,10] = 2.4
$row4[0,11] = 6.3
$row4[0,12] = 1.02
$row4[0,13] = 15
$row4[0,14] = 1.15
$row4[0,15] = 4.05
$row4[0,16] = 1.6
$row4[0,17] = 2.07
$row4[0,18] = 1.29
$row4[0,19] = 3.32
$row4[0,20] = 1.91
$row4[0,21] = 1.85
$row4[0,22] = 6.4
$row4[0,23] = 5.7
$row4[0,24] = 7.2
$row4[0,25] = 7.3
$row4[0,26] = 9.25
$row4[0,27] = 21
$row4[0,28] = 13
$row4[0,29] = 8
$row4[0,30] = 16.5
$row4[0,31] = 65
$row4[0,32] = 450
$row4[0,33] = 16
$row4[0,34] = 35
$row4[0,35] = 18
$row4[0,36] = 110
$row4[0,37] = 60
$row4[0,38] = 50
$row4[0,39] = 3.15
$row4[0,40] = 6
$row4[0,41] = 16
$row4[0,42] = 16
$row4[0,43] = 45
$row4[0,44] = 200
$row4[0,45] = 3.1
$row4[0,46] = 8.5
$row4[0,47] = 80
$row4[0,48] = 8.25
$row4[0,49] = 40
$row4[0,50] = 40
$row4[0,51] = 300
$row4[0,52] = 300
$row4[0,53] = 400
$row4[0,54] = 51
$row4[0,55] = 51
$ws.Range("A4:BD4").Value = $row4

# New row appended for Eldense vs Huesca (Spain - LaLiga2, 16:00).
$row5 = New-Object 'object[,]' 1,56
$row5[0,0] = "zuhsk28K"
$row5[0,1] = "20/11/2024"
$row5[0,2] = "16:00"
$row5[0,3] = "SPAIN - LALIGA2"
$row5[0,4] = "Eldense"
$row5[0,5] = "Huesca"
$row5[0,6] = 2.35
$row5[0,7] = 2.7
$row5[0,8] = 3.7
$row5[0,9] = 3.4
$row5[0,10] = 1.8
$row5[0,11] = 4.5
$row5[0,12] = 1.17
$row5[0,13] = 5
$row5[0,14] = 1.67
$row5[0,15] = 2.1
$row5[0,16] = 3.4
$row5[0,17] = 1.33
$row5[0,18] = 1.75
$row5[0,19] = 2.05
$row5[0,20] = 2.5
$row5[0,21] = 1.5
$row5[0,22] = 5
$row5[0,23] = 9.5
$row5[0,24] = 11
$row5[0,25] = 23
$row5[0,26] = 29
$row5[0,27] = 51
$row5[0,28] = 4.75
$row5[0,29] = 6
$row5[0,30] = 23
$row5[0,31] = 101
$row5[0,32] = 201
$row5[0,33] = 7
$row5[0,34] = 15
$row5[0,35] = 15
$row5[0,36] = 41
$row5[0,37] = 41
$row5[0,38] = 51
$row5[0,39] = 4
$row5[0,40] = 15
$row5[0,41] = 34
$row5[0,42] = 51
$row5[0,43] = 101
$row5[0,44] = 351
$row5[0,45] = 2
$row5[0,46] = 11
$row5[0,47] = 101
$row5[0,48] = 5
$row5[0,49] = 23
$row5[0,50] = 41
$row5[0,51] = 81
$row5[0,52] = 151
$row5[0,53] = 500
$row5[0,54] = 81
$row5[0,55] = 81
$ws.Range("A5:BD5").Value = $row5
